$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value2 = 812.3333
$ws.Range("J19").Value2 = 891.4286
$ws.Range("L19").Value2 = 891.4286
$ws.Range("N19").Value2 = -1241.4286

$ws.Range("H43").Value2 = 1947
$ws.Range("I43").Value2 = 1947
$ws.Range("K43").Value2 = 1947
$ws.Range("M43").Value2 = -1878

$ws.Range("H70").Value2 = 0
$ws.Range("J70").Value2 = 0
$ws.Range("L70").Value2 = 0
$ws.Range("N70").Value2 = $null

$ws.Range("H73").Value2 = 0
$ws.Range("J73").Value2 = 0
$ws.Range("L73").Value2 = 0
$ws.Range("N73").Value2 = $null

$ws.Range("H86").Value2 = 8263.478999999999
$ws.Range("I86").Value2 = 7852.2856
$ws.Range("K86").Value2 = 7852.2856
$ws.Range("M86").Value2 = -6729.2856

$ws.Range("H89").Value2 = 8263.478999999999
$ws.Range("I89").Value2 = 7852.2856
$ws.Range("K89").Value2 = 39261.428
$ws.Range("M89").Value2 = -33645.428

$ws.Range("H92").Value2 = 1728.9286
$ws.Range("J92").Value2 = 300
$ws.Range("L92").Value2 = 300
$ws.Range("N92").Value2 = -2796

$ws.Range("H107").Value2 = 2115.7856
$ws.Range("I107").Value2 = 2354.24
$ws.Range("K107").Value2 = 2354.24
$ws.Range("M107").Value2 = -434.2399999999998

$ws.Range("H112").Value2 = 1878.1
$ws.Range("J112").Value2 = 1640.8572
$ws.Range("L112").Value2 = 4922.571599999999
$ws.Range("N112").Value2 = -7138.571599999999

$ws.Range("H113").Value2 = 3428.64
$ws.Range("I113").Value2 = 3411.1333
$ws.Range("J113").Value2 = 3454.9
$ws.Range("K113").Value2 = 3411.1333
$ws.Range("L113").Value2 = 3454.9
$ws.Range("M113").Value2 = -157.1333
$ws.Range("N113").Value2 = -9962.9

$ws.Range("H125").Value2 = 4352.4614
$ws.Range("I125").Value2 = 4358.4
$ws.Range("K125").Value2 = 39225.6
$ws.Range("M125").Value2 = -36765.6

$ws.Range("H132").Value2 = 13856.294
$ws.Range("I132").Value2 = 16889
$ws.Range("K132").Value2 = 50667
$ws.Range("M132").Value2 = -48137

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value2 = 26616.5
$ws.Range("J44").Value2 = 26616.5
$ws.Range("L44").Value2 = 26616.5
$ws.Range("N44").Value2 = -27592.5

$ws.Range("H55").Value2 = 29400
$ws.Range("J55").Value2 = 29400
$ws.Range("L55").Value2 = 29400
$ws.Range("N55").Value2 = -30030

$ws.Range("H74").Value2 = 2659782.2
$ws.Range("I74").Value2 = 3277848
$ws.Range("K74").Value2 = 3277848
$ws.Range("M74").Value2 = -3276974

$ws.Range("H77").Value2 = 2659782.2
$ws.Range("I77").Value2 = 3277848
$ws.Range("K77").Value2 = 16389240
$ws.Range("M77").Value2 = -16384872

$ws.Range("H94").Value2 = 0
$ws.Range("J94").Value2 = 0
$ws.Range("L94").Value2 = 0
$ws.Range("N94").Value2 = $null

$ws.Range("H97").Value2 = 33334176
$ws.Range("J97").Value2 = 111111704
$ws.Range("L97").Value2 = 111111704
$ws.Range("N97").Value2 = -111112696

$ws.Range("H132").Value2 = 6503.0415
$ws.Range("I132").Value2 = 4597.6665
$ws.Range("J132").Value2 = 7646.2666
$ws.Range("K132").Value2 = 13792.9995
$ws.Range("L132").Value2 = 22938.7998
$ws.Range("M132").Value2 = -11262.9995
$ws.Range("N132").Value2 = -27998.7998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value2 = 1227.2
$ws.Range("I37").Value2 = 1259.625
$ws.Range("K37").Value2 = 1259.625
$ws.Range("M37").Value2 = -1122.625

$ws.Range("H126").Value2 = 30000
$ws.Range("J126").Value2 = 30000
$ws.Range("L126").Value2 = 30000
$ws.Range("N126").Value2 = -39880

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value2 = 88956.53999999999
$ws.Range("J16").Value2 = 4341.6665
$ws.Range("L16").Value2 = 4341.6665
$ws.Range("N16").Value2 = -4915.6665

$ws.Range("H58").Value2 = 10253508
$ws.Range("I58").Value2 = 20835526
$ws.Range("J58").Value2 = 4206639.5
$ws.Range("K58").Value2 = 20835526
$ws.Range("L58").Value2 = 4206639.5
$ws.Range("M58").Value2 = -20835323
$ws.Range("N58").Value2 = -4207045.5

$ws.Range("H68").Value2 = 41832.668
$ws.Range("J68").Value2 = 41832.668
$ws.Range("L68").Value2 = 41832.668
$ws.Range("N68").Value2 = -43330.668

$ws.Range("H71").Value2 = 41832.668
$ws.Range("J71").Value2 = 41832.668
$ws.Range("L71").Value2 = 125498.004
$ws.Range("N71").Value2 = -132986.004

$ws.Range("H86").Value2 = 32503.428
$ws.Range("I86").Value2 = 62446.07
$ws.Range("J86").Value2 = 12541.667
$ws.Range("K86").Value2 = 62446.07
$ws.Range("L86").Value2 = 12541.667
$ws.Range("M86").Value2 = -61323.07
$ws.Range("N86").Value2 = -14787.667

$ws.Range("H89").Value2 = 32503.428
$ws.Range("I89").Value2 = 62446.07
$ws.Range("J89").Value2 = 12541.667
$ws.Range("K89").Value2 = 312230.35
$ws.Range("L89").Value2 = 62708.335
$ws.Range("M89").Value2 = -306614.35
$ws.Range("N89").Value2 = -73940.33499999999

$ws.Range("H113").Value2 = 88956.53999999999
$ws.Range("J113").Value2 = 4341.6665
$ws.Range("L113").Value2 = 4341.6665
$ws.Range("N113").Value2 = -8681.666499999999

$ws.Range("H136").Value2 = 10253508
$ws.Range("I136").Value2 = 20835526
$ws.Range("J136").Value2 = 4206639.5
$ws.Range("K136").Value2 = 62506578
$ws.Range("L136").Value2 = 12619918.5
$ws.Range("M136").Value2 = -62504028
$ws.Range("N136").Value2 = -12625018.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value2 = 909.1539
$ws.Range("I50").Value2 = 513.1818
$ws.Range("J50").Value2 = 3087
$ws.Range("K50").Value2 = 1539.5454
$ws.Range("L50").Value2 = 9261
$ws.Range("M50").Value2 = -1058.5454
$ws.Range("N50").Value2 = -10223

$ws.Range("H53").Value2 = 909.1539
$ws.Range("I53").Value2 = 513.1818
$ws.Range("J53").Value2 = 3087
$ws.Range("K53").Value2 = 1539.5454
$ws.Range("L53").Value2 = 9261
$ws.Range("M53").Value2 = -1058.5454
$ws.Range("N53").Value2 = -10223

$ws.Range("H68").Value2 = 5769
$ws.Range("J68").Value2 = 5955.3335
$ws.Range("L68").Value2 = 17866.0005
$ws.Range("N68").Value2 = -19488.0005

$ws.Range("H71").Value2 = 5769
$ws.Range("J71").Value2 = 5955.3335
$ws.Range("L71").Value2 = 53598.0015
$ws.Range("N71").Value2 = -61710.0015

$ws.Range("H107").Value2 = 3502.0625
$ws.Range("J107").Value2 = 4494.909
$ws.Range("L107").Value2 = 13484.727
$ws.Range("N107").Value2 = -17324.727

$ws.Range("H117").Value2 = 12077.333
$ws.Range("I117").Value2 = 154.66667
$ws.Range("K117").Value2 = 464.00001
$ws.Range("M117").Value2 = 2977.99999

$ws.Range("H137").Value2 = 6249
$ws.Range("I137").Value2 = 2991.8
$ws.Range("J137").Value2 = 11677.667
$ws.Range("K137").Value2 = 8975.400000000001
$ws.Range("L137").Value2 = 35033.001
$ws.Range("M137").Value2 = -3875.400000000001
$ws.Range("N137").Value2 = -45233.001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 11501.3545
$ws.Range("I70").Value2 = 13622.954
$ws.Range("J70").Value2 = 6315.222
$ws.Range("K70").Value2 = 13622.954
$ws.Range("L70").Value2 = 6315.222
$ws.Range("M70").Value2 = -13352.954
$ws.Range("N70").Value2 = -6855.222

$ws.Range("H73").Value2 = 11501.3545
$ws.Range("I73").Value2 = 13622.954
$ws.Range("J73").Value2 = 6315.222
$ws.Range("K73").Value2 = 13622.954
$ws.Range("L73").Value2 = 6315.222
$ws.Range("M73").Value2 = -12686.954
$ws.Range("N73").Value2 = -8187.222

$ws.Range("H102").Value2 = 2801.8572
$ws.Range("I102").Value2 = 977.1667
$ws.Range("J102").Value2 = 13750
$ws.Range("K102").Value2 = 977.1667
$ws.Range("L102").Value2 = 13750
$ws.Range("M102").Value2 = 644.8333
$ws.Range("N102").Value2 = -16994

$ws.Range("H122").Value2 = 40689
$ws.Range("I122").Value2 = 51795.3
$ws.Range("J122").Value2 = 3668
$ws.Range("K122").Value2 = 155385.9
$ws.Range("L122").Value2 = 11004
$ws.Range("M122").Value2 = -152935.9
$ws.Range("N122").Value2 = -15904

$ws.Range("H132").Value2 = 14899.9375
$ws.Range("I132").Value2 = 10673.75
$ws.Range("K132").Value2 = 32021.25
$ws.Range("M132").Value2 = -29491.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value2 = 8887
$ws.Range("I61").Value2 = 10190.909
$ws.Range("J61").Value2 = 5301.25
$ws.Range("K61").Value2 = 10190.909
$ws.Range("L61").Value2 = 5301.25
$ws.Range("M61").Value2 = -9988.909
$ws.Range("N61").Value2 = -5705.25

$ws.Range("H100").Value2 = 4299.75
$ws.Range("J100").Value2 = 4299.75
$ws.Range("L100").Value2 = 4299.75
$ws.Range("N100").Value2 = -5381.75

$ws.Range("H113").Value2 = 8887
$ws.Range("I113").Value2 = 10190.909
$ws.Range("J113").Value2 = 5301.25
$ws.Range("K113").Value2 = 10190.909
$ws.Range("L113").Value2 = 5301.25
$ws.Range("M113").Value2 = -8020.909
$ws.Range("N113").Value2 = -9641.25

$ws.Range("H118").Value2 = 0
$ws.Range("J118").Value2 = 0
$ws.Range("L118").Value2 = 0
$ws.Range("N118").Value2 = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value2 = 64990
$ws.Range("J92").Value2 = 64990
$ws.Range("L92").Value2 = 64990
$ws.Range("N92").Value2 = -69982

$ws.Range("H132").Value2 = 5210088
$ws.Range("I132").Value2 = 6174390
$ws.Range("K132").Value2 = 18523170
$ws.Range("M132").Value2 = -18520640

$ws.Range("H136").Value2 = 11776751
$ws.Range("I136").Value2 = 5436344.5
$ws.Range("J136").Value2 = 62500000
$ws.Range("K136").Value2 = 16309033.5
$ws.Range("L136").Value2 = 187500000
$ws.Range("M136").Value2 = -16306483.5
$ws.Range("N136").Value2 = -187505100
